$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in the "Adjusted count" (column B) values for the last Commit block
# (rows 163-172), which were previously blank, causing all the dependent
# formulas in columns D and F to evaluate to 0.
$ws.Range("B163").Value = 2615
$ws.Range("B164").Value = 26
$ws.Range("B165").Value = 762
$ws.Range("B166").Value = 617
$ws.Range("B167").Value = 414
$ws.Range("B171").Value = 599
$ws.Range("B172").Value = 110

# Update the view's scroll position / selection to match where the user
# ended up after entering the data (bottom of the sheet).
$ws.Activate()
$ws.Range("A139").Select()
$excel.ActiveWindow.ScrollRow = 139
$ws.Range("B173").Select()
